$d = $word.ActiveDocument

# Locate the paragraph containing the target text and remove the entire
# paragraph (including its trailing paragraph mark) so the list collapses
# cleanly, matching the diff which deletes the whole <w:p> element.
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.Trim("`r", "`a", "`n", " ")
    if ($text -eq "Add a song to the team playlist") {
        $p.Range.Delete()
        break
    }
}
